# Updated cryptos list on Wed Sep 25 19:49:59 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.271.32'
$ws.Cells.Item(2, 5).Value = '  -1.03%  '
$ws.Cells.Item(3, 4).Value = '2.577.42'
$ws.Cells.Item(3, 5).Value = '  -2.45%  '
$ws.Cells.Item(4, 5).Value = '  +0.03%  '
$ws.Cells.Item(5, 4).Value = '''588.14'
$ws.Cells.Item(5, 5).Value = '  -3.31%  '
$ws.Cells.Item(6, 4).Value = '''150.27'
$ws.Cells.Item(6, 5).Value = '  +1.92%  '
$ws.Cells.Item(7, 5).Value = '  +0.04%  '
$ws.Cells.Item(8, 4).Value = '''0.585'
$ws.Cells.Item(8, 5).Value = '  -0.69%  '
$ws.Cells.Item(9, 5).Value = '  +0.30%  '
$ws.Cells.Item(10, 4).Value = '''5.70'
$ws.Cells.Item(10, 5).Value = '  +1.89%  '
$ws.Cells.Item(11, 5).Value = '  -0.28%  '
$ws.Cells.Item(12, 5).Value = '  -0.49%  '
$ws.Cells.Item(13, 4).Value = '''27.59'
$ws.Cells.Item(13, 5).Value = '  +0.24%  '
$ws.Cells.Item(14, 4).Value = '3.041.67'
$ws.Cells.Item(14, 5).Value = '  -2.34%  '
$ws.Cells.Item(15, 4).Value = '63.065.07'
$ws.Cells.Item(15, 5).Value = '  -1.09%  '
$ws.Cells.Item(16, 5).Value = '  +5.01%  '
$ws.Cells.Item(17, 4).Value = '2.562.24'
$ws.Cells.Item(17, 5).Value = '  -2.45%  '
$ws.Cells.Item(18, 4).Value = '''12.16'
$ws.Cells.Item(18, 5).Value = '  +3.21%  '
$ws.Cells.Item(19, 4).Value = '''4.72'
$ws.Cells.Item(19, 5).Value = '  +3.09%  '
$ws.Cells.Item(20, 4).Value = '''343.86'
$ws.Cells.Item(20, 5).Value = '  -0.92%  '
$ws.Cells.Item(21, 4).Value = '''6.83'
$ws.Cells.Item(21, 5).Value = '  -1.30%  '
$ws.Cells.Item(22, 5).Value = '  -0.09%  '
$ws.Cells.Item(23, 4).Value = '''67.20'
$ws.Cells.Item(23, 5).Value = '  +1.31%  '
$ws.Cells.Item(24, 5).Value = '  +0.85%  '
$ws.Cells.Item(25, 4).Value = '''9.17'
$ws.Cells.Item(25, 5).Value = '  -1.39%  '
$ws.Cells.Item(26, 4).Value = '''1.66'
$ws.Cells.Item(26, 5).Value = '  -1.44%  '
$ws.Cells.Item(27, 4).Value = '''558.21'
$ws.Cells.Item(27, 5).Value = '  -0.79%  '
$ws.Cells.Item(28, 4).Value = '''8.04'
$ws.Cells.Item(28, 5).Value = '  -1.28%  '
$ws.Cells.Item(29, 5).Value = '  +0.79%  '
$ws.Cells.Item(30, 5).Value = '  +0.08%  '
$ws.Cells.Item(31, 4).Value = '''2.02'
$ws.Cells.Item(31, 5).Value = '  -1.51%  '
$ws.Cells.Item(32, 4).Value = '0.0₃0854'
$ws.Cells.Item(32, 5).Value = '  -0.16%  '
$ws.Cells.Item(33, 5).Value = '  -1.39%  '
$ws.Cells.Item(34, 4).Value = '''5.21'
$ws.Cells.Item(34, 5).Value = '  -1.63%  '
$ws.Cells.Item(35, 4).Value = '''166.46'
$ws.Cells.Item(35, 5).Value = '  -1.63%  '
$ws.Cells.Item(36, 5).Value = '  +1.44%  '
$ws.Cells.Item(37, 4).Value = '''0.999'
$ws.Cells.Item(37, 5).Value = '  +0.04%  '
$ws.Cells.Item(38, 5).Value = '  +1.59%  '
$ws.Cells.Item(39, 4).Value = '''1.91'
$ws.Cells.Item(39, 5).Value = '  -1.83%  '
$ws.Cells.Item(40, 5).Value = '  -0.07%  '
$ws.Cells.Item(41, 4).Value = '''165.81'
$ws.Cells.Item(41, 5).Value = '  +0.38%  '
$ws.Cells.Item(42, 4).Value = '''39.59'
$ws.Cells.Item(42, 5).Value = '  -1.17%  '
$ws.Cells.Item(43, 4).Value = '''3.98'
$ws.Cells.Item(43, 5).Value = '  +4.63%  '
$ws.Cells.Item(44, 4).Value = '''22.77'
$ws.Cells.Item(44, 5).Value = '  +3.50%  '
$ws.Cells.Item(45, 5).Value = '  +2.14%  '
$ws.Cells.Item(46, 5).Value = '  +5.76%  '
$ws.Cells.Item(47, 4).Value = '''0.627'
$ws.Cells.Item(47, 5).Value = '  -0.21%  '
$ws.Cells.Item(48, 4).Value = '''0.0250'
$ws.Cells.Item(48, 5).Value = '  +1.61%  '
$ws.Cells.Item(49, 4).Value = '''0.0959'
$ws.Cells.Item(49, 5).Value = '  +0.00%  '
$ws.Cells.Item(50, 4).Value = '''19.17'
$ws.Cells.Item(50, 5).Value = '  +1.60%  '
$ws.Cells.Item(51, 4).Value = '0.0₆0235'
$ws.Cells.Item(51, 5).Value = '  +19.54%  '
